# Update gh-pages to output generated at 456a3b4
# Sheets "展览" and "全部类型" each have F2 (想去人数) and F4 (想去人数) updated.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1443
    $ws.Range("F4").Value = 11
}
